# Scheduled-runner style refresh of the Leve profit tables in Carbuncle_Profits.xlsx.
# For each leve row below, re-applies the latest fetched market-board prices
# (currentAveragePrice / NQ / HQ columns H-L) and the resulting profit figures
# (LeveProfitNQ / LeveProfitHQ columns M-N) across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1113861.2
$ws.Range("J17").Value = 1113861.2
$ws.Range("L17").Value = 3341583.6
$ws.Range("N17").Value = -3341919.6

# Row 40
$ws.Range("H40").Value = 1193.1333
$ws.Range("I40").Value = 1349.25
$ws.Range("J40").Value = 1136.3636
$ws.Range("K40").Value = 1349.25
$ws.Range("L40").Value = 1136.3636
$ws.Range("M40").Value = -1174.25
$ws.Range("N40").Value = -1486.3636

# Row 64
$ws.Range("H64").Value = 3271.92
$ws.Range("I64").Value = 2970.075
$ws.Range("J64").Value = 4479.3
$ws.Range("K64").Value = 2970.075
$ws.Range("L64").Value = 4479.3
$ws.Range("M64").Value = -2722.075
$ws.Range("N64").Value = -4975.3

# Row 67
$ws.Range("H67").Value = 3271.92
$ws.Range("I67").Value = 2970.075
$ws.Range("J67").Value = 4479.3
$ws.Range("K67").Value = 2970.075
$ws.Range("L67").Value = 4479.3
$ws.Range("M67").Value = -2112.075
$ws.Range("N67").Value = -6195.3

# Row 98
$ws.Range("H98").Value = 955.0909
$ws.Range("I98").Value = 821.10345
$ws.Range("J98").Value = 1926.5
$ws.Range("K98").Value = 821.10345
$ws.Range("L98").Value = 1926.5
$ws.Range("M98").Value = 676.89655
$ws.Range("N98").Value = -4922.5

# Row 122
$ws.Range("H122").Value = 955.0909
$ws.Range("I122").Value = 821.10345
$ws.Range("J122").Value = 1926.5
$ws.Range("K122").Value = 2463.31035
$ws.Range("L122").Value = 5779.5
$ws.Range("M122").Value = -13.31034999999974
$ws.Range("N122").Value = -10679.5

# Row 132
$ws.Range("H132").Value = 44780
$ws.Range("I132").Value = 50914.953
$ws.Range("J132").Value = 1835.3334
$ws.Range("K132").Value = 152744.859
$ws.Range("L132").Value = 5506.0002
$ws.Range("M132").Value = -150214.859
$ws.Range("N132").Value = -10566.0002


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4119
$ws.Range("I32").Value = 3543.289
$ws.Range("K32").Value = 3543.289
$ws.Range("M32").Value = -3256.289

# Row 74
$ws.Range("H74").Value = 3864.2683
$ws.Range("I74").Value = 4270.6
$ws.Range("J74").Value = 1494
$ws.Range("K74").Value = 4270.6
$ws.Range("L74").Value = 1494
$ws.Range("M74").Value = -3396.6
$ws.Range("N74").Value = -3242

# Row 77
$ws.Range("H77").Value = 3864.2683
$ws.Range("I77").Value = 4270.6
$ws.Range("J77").Value = 1494
$ws.Range("K77").Value = 21353
$ws.Range("L77").Value = 7470
$ws.Range("M77").Value = -16985
$ws.Range("N77").Value = -16206

# Row 132
$ws.Range("H132").Value = 1544.2174
$ws.Range("I132").Value = 853.7059
$ws.Range("J132").Value = 3500.6667
$ws.Range("K132").Value = 2561.1177
$ws.Range("L132").Value = 10502.0001
$ws.Range("M132").Value = -31.11770000000024
$ws.Range("N132").Value = -15562.0001


$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1282.6154
$ws.Range("I134").Value = 885.4231
$ws.Range("J134").Value = 2871.3845
$ws.Range("K134").Value = 2656.2693
$ws.Range("L134").Value = 8614.1535
$ws.Range("M134").Value = -121.2692999999999
$ws.Range("N134").Value = -13684.1535


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1095.1666
$ws.Range("I16").Value = 1114.2
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1114.2
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -827.2
$ws.Range("N16").Value = -1574

# Row 31
$ws.Range("H31").Value = 3415.7585
$ws.Range("I31").Value = 3053.8
$ws.Range("J31").Value = 3606.2632
$ws.Range("K31").Value = 3053.8
$ws.Range("L31").Value = 3606.2632
$ws.Range("M31").Value = -2758.8
$ws.Range("N31").Value = -4196.263199999999

# Row 34
$ws.Range("H34").Value = 3415.7585
$ws.Range("I34").Value = 3053.8
$ws.Range("J34").Value = 3606.2632
$ws.Range("K34").Value = 3053.8
$ws.Range("L34").Value = 3606.2632
$ws.Range("M34").Value = -2851.8
$ws.Range("N34").Value = -4010.2632

# Row 58
$ws.Range("H58").Value = 1161.1449
$ws.Range("I58").Value = 766.7778
$ws.Range("J58").Value = 2580.8667
$ws.Range("K58").Value = 766.7778
$ws.Range("L58").Value = 2580.8667
$ws.Range("M58").Value = -563.7778
$ws.Range("N58").Value = -2986.8667

# Row 113
$ws.Range("H113").Value = 1095.1666
$ws.Range("I113").Value = 1114.2
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1114.2
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1055.8
$ws.Range("N113").Value = -5340

# Row 134
$ws.Range("H134").Value = 1398.5538
$ws.Range("I134").Value = 1229.3771
$ws.Range("J134").Value = 3978.5
$ws.Range("K134").Value = 3688.1313
$ws.Range("L134").Value = 11935.5
$ws.Range("M134").Value = -1153.1313
$ws.Range("N134").Value = -17005.5

# Row 136
$ws.Range("H136").Value = 1161.1449
$ws.Range("I136").Value = 766.7778
$ws.Range("J136").Value = 2580.8667
$ws.Range("K136").Value = 2300.3334
$ws.Range("L136").Value = 7742.6001
$ws.Range("M136").Value = 249.6666
$ws.Range("N136").Value = -12842.6001


$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 44.846153
$ws.Range("I12").Value = 3.5
$ws.Range("K12").Value = 10.5
$ws.Range("M12").Value = 162.5

# Row 34
$ws.Range("H34").Value = 490.25
$ws.Range("I34").Value = 469.33334
$ws.Range("K34").Value = 1408.00002
$ws.Range("M34").Value = -1324.00002

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

# Row 122
$ws.Range("H122").Value = 1101.2
$ws.Range("I122").Value = 402.09525
$ws.Range("K122").Value = 3618.85725
$ws.Range("M122").Value = -1168.85725

# Row 131
$ws.Range("H131").Value = 5948.1816
$ws.Range("I131").Value = 820
$ws.Range("K131").Value = 2460
$ws.Range("M131").Value = 2580


$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1990.3334
$ws.Range("I132").Value = 1691.9259
$ws.Range("K132").Value = 5075.7777
$ws.Range("M132").Value = -2545.7777


$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1276
$ws.Range("I68").Value = 1060
$ws.Range("J68").Value = 1492
$ws.Range("K68").Value = 1060
$ws.Range("L68").Value = 1492
$ws.Range("M68").Value = -311
$ws.Range("N68").Value = -2990

# Row 71
$ws.Range("H71").Value = 1276
$ws.Range("I71").Value = 1060
$ws.Range("J71").Value = 1492
$ws.Range("K71").Value = 5300
$ws.Range("L71").Value = 7460
$ws.Range("M71").Value = -1556
$ws.Range("N71").Value = -14948

# Row 132
$ws.Range("H132").Value = 5330.8354
$ws.Range("I132").Value = 4951.361
$ws.Range("J132").Value = 7259.8335
$ws.Range("K132").Value = 14854.083
$ws.Range("L132").Value = 21779.5005
$ws.Range("M132").Value = -12324.083
$ws.Range("N132").Value = -26839.5005

# Row 136
$ws.Range("H136").Value = 8131855
$ws.Range("I136").Value = 1788.9117
$ws.Range("K136").Value = 5366.7351
$ws.Range("M136").Value = -2816.7351


$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1389.7433
$ws.Range("I132").Value = 1172.6912
$ws.Range("J132").Value = 3849.6667
$ws.Range("K132").Value = 3518.0736
$ws.Range("L132").Value = 11549.0001
$ws.Range("M132").Value = -988.0735999999997
$ws.Range("N132").Value = -16609.0001

# Row 138
$ws.Range("H138").Value = 50537.25
$ws.Range("J138").Value = 50537.25
$ws.Range("L138").Value = 50537.25
$ws.Range("N138").Value = -60817.25

Write-Output "Updated leve profit rows across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
